$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NATMI ligand-receptor pair table gains a new sending/target cluster
# ("ECs") alongside the existing "FAPs" and "sCs" clusters, producing a
# full 3x3 cross join of cluster pairs (rows 2-10) with refreshed metrics.

$data = @(
    @("ECs", "Pdgfd", "Pdgfrb", "ECs", 2, 0.6666666666666666, 15.78811133333333, 47.364334, 0.7751682960432016, 0.7751682960432017, 2, 0.6666666666666666, 7.259429666666667, 21.778289, 0.05296410708422108, 0.05296410708422109, 114.6126837938362, 1031.514154144526, 0.04105609663992532, 0.04105609663992533),
    @("ECs", "Pdgfd", "Pdgfrb", "FAPs", 2, 0.6666666666666666, 15.78811133333333, 47.364334, 0.7751682960432016, 0.7751682960432017, 3, 1, 60.13240533333334, 180.397216, 0.438720299189682, 0.4387202991896821, 949.3771101437939, 8544.393991294144, 0.3400820667624294, 0.3400820667624295),
    @("ECs", "Pdgfd", "Pdgfrb", "sCs", 2, 0.6666666666666666, 15.78811133333333, 47.364334, 0.7751682960432016, 0.7751682960432017, 3, 1, 69.67135866666666, 209.014076, 0.5083155937260968, 0.5083155937260969, 1099.979167373931, 9899.812506365384, 0.3940301326408468, 0.3940301326408469),
    @("FAPs", "Pdgfd", "Pdgfrb", "ECs", 3, 1, 3.325805333333333, 9.977416, 0.1632911498266644, 0.1632911498266645, 2, 0.6666666666666666, 7.259429666666667, 21.778289, 0.05296410708422108, 0.05296410708422109, 24.14344990235822, 217.291049121224, 0.008648569945325044, 0.008648569945325047),
    @("FAPs", "Pdgfd", "Pdgfrb", "FAPs", 3, 1, 3.325805333333333, 9.977416, 0.1632911498266644, 0.1632911498266645, 3, 1, 60.13240533333334, 180.397216, 0.438720299189682, 0.4387202991896821, 199.9886743637618, 1799.898069273856, 0.07163914210698141, 0.07163914210698144),
    @("FAPs", "Pdgfd", "Pdgfrb", "sCs", 3, 1, 3.325805333333333, 9.977416, 0.1632911498266644, 0.1632911498266645, 3, 1, 69.67135866666666, 209.014076, 0.5083155937260968, 0.5083155937260969, 231.7133762341795, 2085.420386107616, 0.08300343777435795, 0.08300343777435799),
    @("sCs", "Pdgfd", "Pdgfrb", "ECs", 3, 1, 1.253417, 3.760251, 0.06154055413013397, 0.06154055413013397, 2, 0.6666666666666666, 7.259429666666667, 21.778289, 0.05296410708422108, 0.05296410708422109, 9.099092554504333, 81.891832990539, 0.00325944049897072, 0.00325944049897072),
    @("sCs", "Pdgfd", "Pdgfrb", "FAPs", 3, 1, 1.253417, 3.760251, 0.06154055413013397, 0.06154055413013397, 3, 1, 60.13240533333334, 180.397216, 0.438720299189682, 0.4387202991896821, 75.37097909569067, 678.3388118612161, 0.0269990903202712, 0.0269990903202712),
    @("sCs", "Pdgfd", "Pdgfrb", "sCs", 3, 1, 1.253417, 3.760251, 0.06154055413013397, 0.06154055413013397, 3, 1, 69.67135866666666, 209.014076, 0.5083155937260968, 0.5083155937260969, 87.32726536589733, 785.9453882930759, 0.03128202331089205, 0.03128202331089205)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$col - 1]
    }
}
